$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (NumberFormat "@") for Price cells whose new value
# would otherwise be auto-parsed as a number by Excel, losing the
# original fixed-decimal text formatting used throughout column D.
$textForceCells = @("D4", "D5", "D6", "D8", "D10", "D11", "D13", "D15", "D19", "D21", "D22", "D24", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D37", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '69.360.50'
$ws.Range("E2").Value = '  -2.28%  '
$ws.Range("D3").Value = '3.682.97'
$ws.Range("E3").Value = '  -2.87%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '684.82'
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("D6").Value = '159.58'
$ws.Range("E6").Value = '  -5.86%  '
$ws.Range("D7").Value = '3.681.24'
$ws.Range("E7").Value = '  -2.84%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("E9").Value = '  -5.58%  '
$ws.Range("D10").Value = '0.145'
$ws.Range("E10").Value = '  -8.73%  '
$ws.Range("D11").Value = '7.12'
$ws.Range("E11").Value = '  -5.39%  '
$ws.Range("E12").Value = '  -8.69%  '
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").Value = '  -6.59%  '
$ws.Range("D14").Value = '4.305.70'
$ws.Range("E14").Value = '  -2.94%  '
$ws.Range("D15").Value = '32.40'
$ws.Range("E15").Value = '  -10.15%  '
$ws.Range("D16").Value = '3.686.35'
$ws.Range("E16").Value = '  -5.09%  '
$ws.Range("D17").Value = '69.376.09'
$ws.Range("E17").Value = '  -2.62%  '
$ws.Range("E18").Value = '  -0.76%  '
$ws.Range("D19").Value = '15.89'
$ws.Range("E19").Value = '  -9.20%  '
$ws.Range("E20").Value = '  -9.82%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '10.18'
$ws.Range("E21").Value = '  -1.84%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '468.08'
$ws.Range("E22").Value = '  -9.04%  '
$ws.Range("E23").Value = '  -9.26%  '
$ws.Range("D24").Value = '79.36'
$ws.Range("E24").Value = '  -5.00%  '
$ws.Range("D25").Value = '3.829.29'
$ws.Range("E25").Value = '  -2.94%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '0.0000123'
$ws.Range("E27").Value = '  -11.38%  '
$ws.Range("D28").Value = '10.94'
$ws.Range("E28").Value = '  -12.73%  '
$ws.Range("D29").Value = '9.19'
$ws.Range("E29").Value = '  -9.85%  '
$ws.Range("D30").Value = '2.69'
$ws.Range("E30").Value = '  -8.70%  '
$ws.Range("D31").Value = '1.72'
$ws.Range("E31").Value = '  -12.47%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '2.01'
$ws.Range("E32").Value = '  -10.48%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '6.61'
$ws.Range("E33").Value = '  -8.98%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  -8.09%  '
$ws.Range("E36").Value = '  -4.84%  '
$ws.Range("D37").Value = '8.16'
$ws.Range("E37").Value = '  -11.78%  '
$ws.Range("E38").Value = '  -6.19%  '
$ws.Range("E39").Value = '  -3.68%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  -9.60%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("E43").Value = '  -6.61%  '
$ws.Range("D44").Value = '166.14'
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("D45").Value = '47.76'
$ws.Range("E45").Value = '  -2.73%  '
$ws.Range("D46").Value = '2.72'
$ws.Range("E46").Value = '  -14.48%  '
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("E48").Value = '  -4.77%  '
$ws.Range("D49").Value = '0.000273'
$ws.Range("E49").Value = '  -8.84%  '
$ws.Range("D50").Value = '28.08'
$ws.Range("E50").Value = '  -4.58%  '
$ws.Range("D51").Value = '7.79'
$ws.Range("E51").Value = '  -9.35%  '
